# Rename the existing sheet "Ark1" -> "Statistikk"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "Statistikk"

# Add the new " BIM" sheet at the end and fill it in
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = " BIM"
$ws2.Range("A1").Value = "LUK:"
$ws2.Range("B1").Value = "LUK:"
$ws2.Range("C1").Value = "LUK:"
$ws2.Range("A2").Value = "BYTS1401"
$ws2.Range("B2").Value = "BYFE1201"
$ws2.Range("A3").Value = "BYFE3100"
$ws2.Range("B3").Value = "BYTS1401"
$ws2.Range("B4").Value = "BYFE1201"
$ws2.Range("B5").Value = "EMPE2500"

# Add the new "modell" sheet at the end and fill it in
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "modell"
$ws3.Range("A1").Value = "LUK:"
$ws3.Range("B1").Value = "LUK:"
$ws3.Range("C1").Value = "LUK:"
$ws3.Range("A2").Value = "BYPE1500"
$ws3.Range("B2").Value = "BYFE1000"
$ws3.Range("C2").Value = "BYFE1201"
$ws3.Range("A3").Value = "BYPE2200"
$ws3.Range("B3").Value = "BYTS1401"
$ws3.Range("C3").Value = "BYPE2700"
$ws3.Range("A4").Value = "BYPE2700"
$ws3.Range("B4").Value = "DAVE3705"
$ws3.Range("C4").Value = "BYFE1201"
$ws3.Range("A5").Value = "BYTS2691"
$ws3.Range("B5").Value = "EMFE1000"
$ws3.Range("C5").Value = "EMTS2200"
$ws3.Range("A6").Value = "BYVE3200"
$ws3.Range("B6").Value = "EMPE2500"
$ws3.Range("C6").Value = "EMTS2300"
$ws3.Range("A7").Value = "FEPE2100"
$ws3.Range("B7").Value = "EMVE3700"
$ws3.Range("C7").Value = "EMVE3700"
$ws3.Range("A8").Value = "EMVE3700"
$ws3.Range("B8").Value = "DAVE3705"

# The two pie charts on Statistikk reference the old sheet name via cached
# SERIES formulas - repoint them at the renamed sheet.
$co1 = $ws1.ChartObjects().Item(1)
$s1 = $co1.Chart.SeriesCollection(1)
$s1.Formula = "=SERIES(,Statistikk!`$G`$1:`$H`$1,Statistikk!`$G`$2:`$H`$2,1)"

$co2 = $ws1.ChartObjects().Item(2)
$s2 = $co2.Chart.SeriesCollection(1)
$s2.Formula = "=SERIES(,Statistikk!`$G`$1:`$H`$1,Statistikk!`$G`$2:`$H`$2,1)"

# Restore Statistikk as the active/selected sheet and update its selection
$ws1.Activate()
$ws1.Range("J27").Select()
